# Update automàtic: dades i banners [2026-02-05 15:49]
# Refresh the DATA_EXTRACCIO (column E) timestamps for rows 2-36 on the
# active sheet to the new extraction times.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = [ordered]@{
    2  = "2026-02-05 15:47:47"
    3  = "2026-02-05 15:47:49"
    4  = "2026-02-05 15:47:52"
    5  = "2026-02-05 15:47:54"
    6  = "2026-02-05 15:47:56"
    7  = "2026-02-05 15:47:59"
    8  = "2026-02-05 15:48:01"
    9  = "2026-02-05 15:48:04"
    10 = "2026-02-05 15:48:06"
    11 = "2026-02-05 15:48:08"
    12 = "2026-02-05 15:48:11"
    13 = "2026-02-05 15:48:13"
    14 = "2026-02-05 15:48:15"
    15 = "2026-02-05 15:48:18"
    16 = "2026-02-05 15:48:20"
    17 = "2026-02-05 15:48:23"
    18 = "2026-02-05 15:48:25"
    19 = "2026-02-05 15:48:28"
    20 = "2026-02-05 15:48:30"
    21 = "2026-02-05 15:48:32"
    22 = "2026-02-05 15:48:35"
    23 = "2026-02-05 15:48:37"
    24 = "2026-02-05 15:48:40"
    25 = "2026-02-05 15:48:42"
    26 = "2026-02-05 15:48:45"
    27 = "2026-02-05 15:48:47"
    28 = "2026-02-05 15:48:50"
    29 = "2026-02-05 15:48:53"
    30 = "2026-02-05 15:48:55"
    31 = "2026-02-05 15:48:58"
    32 = "2026-02-05 15:49:00"
    33 = "2026-02-05 15:49:03"
    34 = "2026-02-05 15:49:05"
    35 = "2026-02-05 15:49:08"
    36 = "2026-02-05 15:49:10"
}

foreach ($row in $timestamps.Keys) {
    $ws.Range("E$row").Value = $timestamps[$row]
}
